$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy formatting from the last existing data row (1024) down to the new rows (1025-1042)
$src = $ws.Range("A1024:C1024")
$dst = $ws.Range("A1025:C1042")
$src.Copy($dst)

# Fill in the new translation rows
$ws.Cells.Item(1025, 1).Value = "cs"
$ws.Cells.Item(1025, 2).Value = 'lab.cotton.label'
$ws.Cells.Item(1025, 3).Value = 'Vaty'

$ws.Cells.Item(1026, 1).Value = "cs"
$ws.Cells.Item(1026, 2).Value = 'lab.cotton.title'
$ws.Cells.Item(1026, 3).Value = 'Vaty'

$ws.Cells.Item(1027, 1).Value = "cs"
$ws.Cells.Item(1027, 2).Value = 'lab.cotton.filter.title'
$ws.Cells.Item(1027, 3).Value = 'Filtr vat'

$ws.Cells.Item(1028, 1).Value = "cs"
$ws.Cells.Item(1028, 2).Value = 'lab.cotton.button.create'
$ws.Cells.Item(1028, 3).Value = 'Nová vata'

$ws.Cells.Item(1029, 1).Value = "cs"
$ws.Cells.Item(1029, 2).Value = 'lab.cotton.table.name'
$ws.Cells.Item(1029, 3).Value = 'Název'

$ws.Cells.Item(1030, 1).Value = "cs"
$ws.Cells.Item(1030, 2).Value = 'lab.cotton.table.vendor'
$ws.Cells.Item(1030, 3).Value = 'Výrobce'

$ws.Cells.Item(1031, 1).Value = "cs"
$ws.Cells.Item(1031, 2).Value = 'lab.cotton.button.create'
$ws.Cells.Item(1031, 3).Value = 'Nová vata'

$ws.Cells.Item(1032, 1).Value = "cs"
$ws.Cells.Item(1032, 2).Value = 'lab.cotton.context.menu'
$ws.Cells.Item(1032, 3).Value = 'Vata [{{data.name}}]'

$ws.Cells.Item(1033, 1).Value = "cs"
$ws.Cells.Item(1033, 2).Value = 'lab.cotton.button.edit'
$ws.Cells.Item(1033, 3).Value = 'Upravit vatu'

$ws.Cells.Item(1034, 1).Value = "cs"
$ws.Cells.Item(1034, 2).Value = 'lab.cotton.button.delete'
$ws.Cells.Item(1034, 3).Value = 'Odstranit vatu'

$ws.Cells.Item(1035, 1).Value = "cs"
$ws.Cells.Item(1035, 2).Value = 'lab.cotton.button.delete.confirm.title'
$ws.Cells.Item(1035, 3).Value = 'Odstranit vatu'

$ws.Cells.Item(1036, 1).Value = "cs"
$ws.Cells.Item(1036, 2).Value = 'lab.cotton.button.delete.confirm'
$ws.Cells.Item(1036, 3).Value = 'Opravdu si přejete odstranit vybranou vatu? Tato akce pravděpodobně smaže velké množství jiných dat, hlavně pak záznamy o vapování. Není cesty zpět, postupujte tedy prosím obezřetně.'

$ws.Cells.Item(1037, 1).Value = "cs"
$ws.Cells.Item(1037, 2).Value = 'lab.cotton.button.delete.confirm.ok'
$ws.Cells.Item(1037, 3).Value = 'Odstranit vatu'

$ws.Cells.Item(1038, 1).Value = "cs"
$ws.Cells.Item(1038, 2).Value = 'lab.cotton.table.description'
$ws.Cells.Item(1038, 3).Value = 'Popis'

$ws.Cells.Item(1039, 1).Value = "cs"
$ws.Cells.Item(1039, 2).Value = 'lab.cotton.index.label'
$ws.Cells.Item(1039, 3).Value = 'Vata [{{data.name}}]'

$ws.Cells.Item(1040, 1).Value = "cs"
$ws.Cells.Item(1040, 2).Value = 'lab.cotton.index.title'
$ws.Cells.Item(1040, 3).Value = 'Detail vaty'

$ws.Cells.Item(1041, 1).Value = "cs"
$ws.Cells.Item(1041, 2).Value = 'lab.cotton.update.submit'
$ws.Cells.Item(1041, 3).Value = 'Aktulizovat vatu'

$ws.Cells.Item(1042, 1).Value = "cs"
$ws.Cells.Item(1042, 2).Value = 'lab.cotton.updated.message'
$ws.Cells.Item(1042, 3).Value = 'Vata [{{data.name}}] byla aktualizována.'

# Update sheet view to match the new scroll/selection position
$ws.Activate()
try {
    $ws.Application.ActiveWindow.ScrollRow = 1018
    $ws.Application.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B1030").Select()